# Applies the "working MVP of screening, looping over family members" edit:
#  - Consolidates several split w:r runs back into single runs (achieved by
#    running Find & Replace over the affected text so Word's save pipeline
#    re-merges adjacent same-formatted runs).
#  - Relocates the stray _GoBack bookmark from the end of the "Shelter
#    letter" checklist item down into the ever_arrested yes/no field,
#    adding the missing ")" before the closing "}}".
#
# NB: Find.Execute signature used throughout:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#           Format, ReplaceWith, Replace)

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# 1. "Children under" / " 14 " -> "Children under 14 "
Replace-Text "Children under 14 " "Children under 14 "

# 2. Drop the old _GoBack bookmark sitting after "Shelter letter"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. "Supe" / "rvisor name: ___________" -> "Supervisor name: ___________"
Replace-Text "Supervisor name: ___________" "Supervisor name: ___________"

# 4. "Applicant eligible" / "—" / "Prepare EAD (Employment Authorization Documents)"
Replace-Text "Applicant eligible—Prepare EAD (Employment Authorization Documents)" "Applicant eligible—Prepare EAD (Employment Authorization Documents)"

# 5. "Referred out" / "—" / "Explain why applicant can" / "’" / "t be seen today. ..."
Replace-Text "Referred out—Explain why applicant can’t be seen today. Ensure applicant takes this folder with legal screener to check-out station. " "Referred out—Explain why applicant can’t be seen today. Ensure applicant takes this folder with legal screener to check-out station. "

# 6. "Indica" / "te any addendum needed to be included if eligible:"
Replace-Text "Indicate any addendum needed to be included if eligible:" "Indicate any addendum needed to be included if eligible:"

# 7. "# of children applying with you" / " today:"
Replace-Text "# of children applying with you today:" "# of children applying with you today:"

# 8. "What " / "is the proof of parole? "
Replace-Text "What is the proof of parole? " "What is the proof of parole? "

# 9. "If no proof --> Flag sup" / "ervisor "
Replace-Text "If no proof --> Flag supervisor " "If no proof --> Flag supervisor "

# 10. " No " / "–" / " " -> " No – "
Replace-Text " No – " " No – "

# 11. " a Final Order of " / "Removal?" -> " a Final Order of Removal?"
#     (narrow match so the preceding standalone "d" run is untouched)
Replace-Text "Removal?" "Removal?"

# 12. " Yes " / "–" / " " -> " Yes – "
Replace-Text " Yes – " " Yes – "

# 13. Insert the missing ")" before " }}" in the ever_arrested field, and
#     plant the _GoBack bookmark right after it.
Replace-Text "ever_arrested }}" "ever_arrested) }}"

$r = $d.Content.Find.Execute("ever_arrested) ")
$bmRange = $d.Content.Find.Parent.Duplicate
$found = $d.Content.Find.Execute("ever_arrested)")
$rng = $d.Content.Duplicate
$rng.Find.Execute("ever_arrested)") | Out-Null
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)

# 14. "Flag to supervisors" / " and " / "list when/where/outcome:"
Replace-Text "and list when/where/outcome:" "and list when/where/outcome:"
